$p = $ppt.ActivePresentation

# 1. Slide 1 notes: clear the numbered "What is inca..." bullet list, leaving an
#    empty trailing bullet paragraph.
$s1 = $p.Slides.Item(1)
$s1.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# 2. Slide 2 notes: clear the long "user-level monitoring" talk-track text.
$s2 = $p.Slides.Item(2)
$s2.NotesPage.Shapes.Item(3).TextFrame.TextRange.Text = ""

# 3. Slide 6 notes: clear the "PingPong reports..." comments.
$s6 = $p.Slides.Item(6)
$s6.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# 4. Slide 7 notes: clear the "Collect Grid benchmark measurements using GrASP..." notes.
$s7 = $p.Slides.Item(7)
$s7.NotesPage.Shapes.Item(3).TextFrame.TextRange.Text = ""

# 5. Remove the final slide ("Inca Information") entirely, along with its notes page.
$p.Slides.Item(9).Delete()
